$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / percentage / non-numeric-looking price cells: set value directly.
# This keeps them as text (Excel cannot coerce them into numbers) and preserves
# exact whitespace / formatting, matching the original inline-string cells.
$ws.Range("D2").Value = "51.761.24"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "2.804.70"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("E6").Value = "  +1.82%  "
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +8.18%  "
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("E13").Value = "  +2.86%  "
$ws.Range("E14").Value = "  +2.40%  "
$ws.Range("D15").Value = "3.246.34"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "2.808.81"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").Value = "51.739.29"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("E19").Value = "  +1.63%  "
$ws.Range("E20").Value = "  +2.90%  "
$ws.Range("E21").Value = "  +3.18%  "
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  -3.43%  "
$ws.Range("E29").Value = "  +11.81%  "
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("E31").Value = "  +3.84%  "
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("E34").Value = "  +8.68%  "
$ws.Range("E35").Value = "  +5.86%  "
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("E39").Value = "  +2.38%  "
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("E42").Value = "  -2.33%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("E46").Value = "  +4.52%  "
$ws.Range("D47").Value = "2.114.44"
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("E48").Value = "  +6.73%  "
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("E50").Value = "  -1.27%  "
$ws.Range("E51").Value = "  +7.70%  "

# Price cells whose new text looks like a plain number (e.g. "1.00", "120.00").
# Excel would silently coerce a bare numeric string into a real number and drop
# formatting like trailing zeros, so force text entry with a leading apostrophe,
# then clear the formatting Excel auto-applies so the cell keeps the default style
# (matching the workbook, where these cells have no explicit style).
$ws.Range("D5").Value = "'354.92"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").Value = "'111.51"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").Value = "'0.556"
$ws.Range("D7").ClearFormats()
$ws.Range("D9").Value = "'0.634"
$ws.Range("D9").ClearFormats()
$ws.Range("D12").Value = "'0.0839"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Value = "'20.04"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").Value = "'7.77"
$ws.Range("D14").ClearFormats()
$ws.Range("D19").Value = "'7.65"
$ws.Range("D19").ClearFormats()
$ws.Range("D23").Value = "'70.56"
$ws.Range("D23").ClearFormats()
$ws.Range("D24").Value = "'268.86"
$ws.Range("D24").ClearFormats()
$ws.Range("D27").Value = "'26.17"
$ws.Range("D27").ClearFormats()
$ws.Range("D29").Value = "'38.87"
$ws.Range("D29").ClearFormats()
$ws.Range("D32").Value = "'52.44"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").Value = "'6.15"
$ws.Range("D33").ClearFormats()
$ws.Range("D35").Value = "'0.0884"
$ws.Range("D35").ClearFormats()
$ws.Range("D36").Value = "'0.0445"
$ws.Range("D36").ClearFormats()
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").Value = "'18.81"
$ws.Range("D38").ClearFormats()
$ws.Range("D43").Value = "'2.22"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").Value = "'120.00"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").Value = "'21.97"
$ws.Range("D45").ClearFormats()
$ws.Range("D46").Value = "'3.41"
$ws.Range("D46").ClearFormats()
$ws.Range("D48").Value = "'2.43"
$ws.Range("D48").ClearFormats()
$ws.Range("D50").Value = "'5.47"
$ws.Range("D50").ClearFormats()
